# Re-write json config, copy xlsx config: add two new source sheets
# ("HEYCAR" and "CINCH") to the VO_UK config workbook, cloned from the
# existing "CARANDCLASSIC" sheet (same header row / column layout / page
# setup), inserted right after it. Also refreshes a few stale sheet
# selections that moved as a side effect of the author's editing session.

$wb = $excel.ActiveWorkbook

# --- Add HEYCAR and CINCH, copied from CARANDCLASSIC (same A1:F1 header
#     row using the shared "name/columns/parameters/logs/field category/
#     type of charts desired" strings, same styles/col widths/page setup)
$src = $wb.Worksheets.Item("CARANDCLASSIC")

$src.Copy($null, $src)
$heycar = $wb.Worksheets.Item($src.Index + 1)
$heycar.Name = "HEYCAR"

$src.Copy($null, $heycar)
$cinch = $wb.Worksheets.Item($heycar.Index + 1)
$cinch.Name = "CINCH"

# --- Refresh the stale selections left over on the sheets touched during
#     the edit session
$vouk = $wb.Worksheets.Item("VO_UK")
$vouk.Activate()
$vouk.Range("C74").Select()

$src.Activate()
$src.Range("A1").Select()

$heycar.Activate()
$heycar.Range("A1").Select()

# CINCH ends up being the last-touched / active sheet
$cinch.Activate()
$cinch.Range("H19").Select()
